# Fruta / hortaliza, semanal
# Insert a new weekly price row at row 18 (pushing existing rows 18-25 down
# to 19-26) and populate it with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 18:25 down one row to make room for the new weekly record.
$ws.Range("A18:R18").EntireRow.Insert()

# Populate the newly inserted row 18 with this week's values.
$ws.Range("A18").Value = 11
$ws.Range("B18").Value = "Vega Monumental Concepción"
$ws.Range("C18").Value = "Bíobío"
$ws.Range("D18").Value = 44468
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = 100112031
$ws.Range("G18").Value = "Poroto verde"
$ws.Range("H18").Value = "Sin especificar"
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 100
$ws.Range("K18").Value = 31000
$ws.Range("L18").Value = 32000
$ws.Range("M18").Value = 31500
$ws.Range("N18").Value = "$/malla 25 kilos"
$ws.Range("O18").Value = "Región de Arica y Parinacota"
$ws.Range("P18").Value = 1260
$ws.Range("Q18").Value = 25
$ws.Range("R18").Value = "Hortaliza"
